$d = $word.ActiveDocument

# 1. "...desirable feature in a game. But creating..." -> "...desirable feature in games. But creating..."
$d.Content.Find.Execute("desirable feature in a game. But creating", $true, $false, $false, $false, $false, $true, 1, $false, "desirable feature in games. But creating", 2)

# 2. "Defective construction of implementation." -> "Defective construction of the implementation."
# (search only the part after the _GoBack bookmark so the bookmark itself is preserved)
$d.Content.Find.Execute("fective construction of implementation.", $true, $false, $false, $false, $false, $true, 1, $false, "fective construction of the implementation.", 2)
